$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 256, shifting existing rows 256:326 down to 257:327
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256 with the new price record
$ws.Cells.Item(256, 1).Value = 4
$ws.Cells.Item(256, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(256, 3).Value = "Los Lagos"
$ws.Cells.Item(256, 4).Value = 44559
$ws.Cells.Item(256, 5).Value = 10
$ws.Cells.Item(256, 6).Value = "Fruta"
$ws.Cells.Item(256, 7).Value = 100102
$ws.Cells.Item(256, 8).Value = "Cítricos"
$ws.Cells.Item(256, 9).Value = 100102005
$ws.Cells.Item(256, 10).Value = "Naranja"
$ws.Cells.Item(256, 11).Value = "Valencia"
$ws.Cells.Item(256, 12).Value = "Primera"
$ws.Cells.Item(256, 13).Value = 200
$ws.Cells.Item(256, 14).Value = 13500
$ws.Cells.Item(256, 15).Value = 14000
$ws.Cells.Item(256, 16).Value = 13750
$ws.Cells.Item(256, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(256, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(256, 19).Value = 917
$ws.Cells.Item(256, 20).Value = 15
